$wb = $excel.ActiveWorkbook
$wsWind = $wb.Worksheets.Item("windspeed")

# --- Update the windspeed sheet: rows 13-76 (Time 11..74) in column B become 15 ---
# This removes the existing RAND()-based formulas (B33:B52) and plain 0/10 values,
# replacing all of them with the constant 15.
for ($r = 13; $r -le 76; $r++) {
    $wsWind.Cells.Item($r, 2).Value2 = 15
}

# --- Add the new "irradiance" worksheet after "windspeed" ---
$wsIrr = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsWind)
$wsIrr.Name = "irradiance"

# Header row
$wsIrr.Cells.Item(1, 1).Value2 = "Time [s]"
$wsIrr.Cells.Item(1, 2).Value2 = "Irradiance [W/m2]"

# Data rows: Time 0..77 in column A, constant Irradiance 1100 in column B
for ($i = 0; $i -le 77; $i++) {
    $row = $i + 2
    $wsIrr.Cells.Item($row, 1).Value2 = $i
    $wsIrr.Cells.Item($row, 2).Value2 = 1100
}

# --- Update selections on both sheets to match the saved view state ---
$wsWind.Activate() | Out-Null
$wsWind.Range("B2:B8").Select() | Out-Null

$wsIrr.Activate() | Out-Null
$wsIrr.Range("B2:B8").Select() | Out-Null
